# "Generate Report for Handback" - fills in the handback columns (Latest
# Target File / Latest Handback File / Latest Handback DateTime) for the
# zh-cn and de-de localization status sheets, flips Status from
# "Ready for handoff" to "Handed back: in sync with en-US", and widens a
# few columns that now hold longer text.
#
# NOTE: Range.Value is not reliable for *reading* in this host (it can
# surface the COM property descriptor instead of the cell's data), so
# reads use Value2/Text; writes still go through Value.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Excel's ColumnWidth is quantized in 1/6-character pixel steps, and the
# stored xlsx <col width> is ColumnWidth + 5/6. Helper to pick the
# ColumnWidth that lands closest to a desired stored width.
function Set-StoredColumnWidth {
    param($col, [double]$storedWidth)
    $col.ColumnWidth = $storedWidth - (5.0/6.0)
}

# ---------------------------------------------------------------------
# Status text: every "Ready for handoff" cell (the zh-cn/de-de "Status"
# column AND the Overview sheet's per-language summary columns, which all
# point at the same shared string) becomes
# "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"
foreach ($addr in @("C2", "C3")) {
    if ($wsZhCn.Range($addr).Value2 -eq $oldStatus) { $wsZhCn.Range($addr).Value = $newStatus }
    if ($wsDeDe.Range($addr).Value2 -eq $oldStatus) { $wsDeDe.Range($addr).Value = $newStatus }
}
foreach ($addr in @("E2", "F2", "E3", "F3")) {
    if ($wsOverview.Range($addr).Value2 -eq $oldStatus) { $wsOverview.Range($addr).Value = $newStatus }
}

# ---------------------------------------------------------------------
# zh-cn / de-de sheets (table columns: ... H=Latest Handoff Datetime,
# I=Latest Target File, J=Latest Handback File, K=Latest Handback DateTime)
# ---------------------------------------------------------------------
$targetFileName = "bba635fa-4500-432b-a42a-59349ebe83e2.md"
$targetUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/125cf44a81c348104c57a0b343aef55eaf147b76/e2e/bba635fa-4500-432b-a42a-59349ebe83e2.md"
$zhHandbackFile = "bba635fa-4500-432b-a42a-59349ebe83e2.23b26500faa09a53e795de58f1ade74bd59d155a.zh-cn.xlf"
$deHandbackFile = "bba635fa-4500-432b-a42a-59349ebe83e2.23b26500faa09a53e795de58f1ade74bd59d155a.de-de.xlf"
$zhHandbackTime = "2016-09-05 21:13:17"
$deHandbackTime = "2016-09-05 21:13:25"

$aDisplay2 = "bba635fa-4500-432b-a42a-59349ebe83e2.md"
$aDisplay3 = "fffffc448ee2-4aa0-444b-8d5a-f4417bd6c085.md"
$aUrl3     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/125cf44a81c348104c57a0b343aef55eaf147b76/e2e/fffffc448ee2-4aa0-444b-8d5a-f4417bd6c085.md"

foreach ($item in @(
        @{ ws = $wsZhCn; handback = $zhHandbackFile; time = $zhHandbackTime },
        @{ ws = $wsDeDe; handback = $deHandbackFile; time = $deHandbackTime }
    )) {
    $ws = $item.ws

    # Latest Target File (I2/I3) + Latest Handback File (J2/J3) + Latest
    # Handback DateTime (K2/K3)
    $ws.Range("I2").Value = $targetFileName
    $ws.Range("J2").Value = $item.handback
    $ws.Range("K2").Value = $item.time

    $ws.Range("I3").Value = $targetFileName
    $ws.Range("J3").Value = $item.handback
    $ws.Range("K3").Value = $item.time

    # Rebuild the hyperlinks in row order (A2, I2, A3, I3) so the new
    # links for the "Latest Target File" cells land right after the
    # matching "Source File Name" link.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $targetUrl, [Type]::Missing, [Type]::Missing, $aDisplay2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I2"), $targetUrl, [Type]::Missing, [Type]::Missing, $aDisplay2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $aUrl3,     [Type]::Missing, [Type]::Missing, $aDisplay3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), $targetUrl, [Type]::Missing, [Type]::Missing, $aDisplay2) | Out-Null

    # Match the blue-underline "HyperLink" look used by the other
    # hyperlink cells in the sheet (A2/A3).
    foreach ($addr in @("I2", "I3")) {
        $f = $ws.Range($addr).Font
        $f.Name = "Calibri"
        $f.Size = 11
        $f.Underline = 2
        $f.Color = 15570276
    }

    # Columns got wider to fit the newly-populated long file names.
    Set-StoredColumnWidth $ws.Columns.Item(3) 29.9777047293527
    Set-StoredColumnWidth $ws.Columns.Item(9) 40
    Set-StoredColumnWidth $ws.Columns.Item(10) 40
}

# ---------------------------------------------------------------------
# Overview sheet: the zh-cn / de-de summary columns (E, F) widen to match.
# ---------------------------------------------------------------------
Set-StoredColumnWidth $wsOverview.Columns.Item(5) 29.9777047293527
Set-StoredColumnWidth $wsOverview.Columns.Item(6) 29.9777047293527
